# Fruta / hortaliza, semanal
# Insert a new weekly record at row 82 of Sheet1, pushing the existing
# rows 82-94 down to 83-95 (new dimension A1:R95).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 82 - everything below shifts down.
$ws.Rows(82).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(82, 1).Value = 7
$ws.Cells.Item(82, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(82, 3).Value = "Ñuble"
$ws.Cells.Item(82, 4).Value = 44918
$ws.Cells.Item(82, 5).Value = 16
$ws.Cells.Item(82, 6).Value = 100112022
$ws.Cells.Item(82, 7).Value = "Arveja Verde"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 60
$ws.Cells.Item(82, 11).Value = 20000
$ws.Cells.Item(82, 12).Value = 21000
$ws.Cells.Item(82, 13).Value = 20500
$ws.Cells.Item(82, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(82, 15).Value = "Región de Ñuble"
$ws.Cells.Item(82, 16).Value = 820
$ws.Cells.Item(82, 17).Value = 25
$ws.Cells.Item(82, 18).Value = "Hortaliza"
